$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# 1. Update the "Description" line.
$ws.Range("A2").Value = "Description: Donor Revenues"

# 2. Update the "Source" line.
$ws.Range("A4").Value = "Source: Local Government Budgets - Ministry of Finance, Planning and Economic Development"

# 3. Insert a new row right after the "Source" line for the new "Source-link" line.
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Source-link: http://www.budget.go.ug/"

# 4. Update the license line (now shifted down to row 16).
$ws.Range("A16").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."

# 5. Insert a new row right after the license line for the licensing-info link.
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
